$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = -8.436
$ws.Range("D6").Value = -8.347
$ws.Range("C7").Value = -13.498
$ws.Range("A10").Value = -20.926
$ws.Range("A12").Value = -21.808
$ws.Range("B13").Value = 6.475
$ws.Range("A18").Value = -21.78
$ws.Range("C20").Value = -13.041
